$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Paragraph 4 ("Esta es otra frase para probar el <commit>"): the word
#    "commit" was wrapped in proofErr spellStart/spellEnd markers and split
#    into its own run. Collapse it back into a single run with no proofing
#    markers, matching the restored/previous revision of the document.
# ---------------------------------------------------------------------------
$p4 = $d.Paragraphs.Item(4)
$xml4 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + `
    '<w:body><w:p w14:paraId="50BBC5AD" w14:textId="13498B87" w:rsidR="009B5A16" w:rsidRDefault="009B5A16">' + `
    '<w:r><w:t>Esta es otra frase para probar el commit</w:t></w:r></w:p></w:body></w:document>' + `
    '</pkg:xmlData></pkg:part></pkg:package>'
$p4.Range.InsertXML($xml4)

# ---------------------------------------------------------------------------
# 2) The "Esta es la frase de prueba creada por Jeiny." paragraph is removed
#    as a standalone paragraph; its (merged) text becomes the run living
#    inside the following paragraph, which only carried an underlined
#    paragraph mark (pPr/rPr/u) before.
# ---------------------------------------------------------------------------
$pUnderline = $d.Paragraphs.Item(6)
$pUnderline.Range.InsertBefore("Esta es la frase de prueba creada por Jeiny.")

$pJeiny = $d.Paragraphs.Item(5)
$pJeiny.Range.Delete()

# ---------------------------------------------------------------------------
# 3) Everything from "Error 1" through the "Errorrrrrrrr ... " paragraph
#    (several blank paragraphs included) is deleted, restoring the document
#    to end right after the underlined paragraph, followed by the final
#    blank paragraph.
# ---------------------------------------------------------------------------
$startPara = $d.Paragraphs.Item(6)
$endPara = $d.Paragraphs.Item(17)
$junkRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$junkRange.Delete()
